# Correction of the reference values (fmod / R columns) on the single
# worksheet. The cells hold numeric-looking text (shared strings), so a
# leading apostrophe is used to force Excel to keep them as text instead
# of re-typing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value  = "'11.320674"
$ws.Range("A4").Value  = "'21.970823"
$ws.Range("A5").Value  = "'24.331106"
$ws.Range("A6").Value  = "'25.77228"
$ws.Range("A7").Value  = "'27.801502"
$ws.Range("A8").Value  = "'29.831356"
$ws.Range("A9").Value  = "'31.347486"
$ws.Range("A10").Value = "'32.429077"
$ws.Range("A11").Value = "'35.616386"
$ws.Range("A12").Value = "'38.407555"
$ws.Range("A13").Value = "'40.773315"
$ws.Range("A14").Value = "'42.723705"
$ws.Range("A15").Value = "'46.29873"
$ws.Range("A16").Value = "'50.558178"
$ws.Range("A17").Value = "'100.85394"

$ws.Range("B3").Value  = "'0.101526074"
$ws.Range("B4").Value  = "'0.26132664"
$ws.Range("B5").Value  = "'0.29006895"
$ws.Range("B6").Value  = "'0.3063978"
$ws.Range("B7").Value  = "'0.32446954"
$ws.Range("B8").Value  = "'0.33473063"
$ws.Range("B9").Value  = "'0.34085318"
$ws.Range("B10").Value = "'0.34529707"
$ws.Range("B11").Value = "'0.3424573"
$ws.Range("B12").Value = "'0.33437043"
$ws.Range("B13").Value = "'0.3230931"
$ws.Range("B14").Value = "'0.310582"
$ws.Range("B15").Value = "'0.28253472"
$ws.Range("B16").Value = "'0.24710885"
$ws.Range("B17").Value = "'0.07500888"
